# "Generate Report for Archive"
# - Update Status value "Ready for handoff" -> "In Translation" (shared across
#   Overview!E2/F2 and the per-locale sheets' Status column, since they all
#   reference the same shared string).
# - Narrow the "Status"-related columns (Overview E:F, and column C on the
#   zh-cn / de-de sheets) to their new width.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newWidth = 13.4101845877511

# Update the status text everywhere it appears (shared string, so setting the
# value on each occurrence keeps everything consistent).
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$zhcn.Range("C2").Value = "In Translation"
$dede.Range("C2").Value = "In Translation"

# Resize the affected columns.
$overview.Columns.Item(5).ColumnWidth = $newWidth
$overview.Columns.Item(6).ColumnWidth = $newWidth
$zhcn.Columns.Item(3).ColumnWidth = $newWidth
$dede.Columns.Item(3).ColumnWidth = $newWidth
